$wb = $excel.ActiveWorkbook

# --- Repayment schedule sheet: insert a new (blank) column before column N ---
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")
$wsSchedule.Columns("N").EntireColumn.Insert()

# Give the newly inserted (blank) column roughly the same width as its
# neighbour to its left (column M), matching Excel's own "insert column"
# behaviour of carrying over the existing column width.
$wsSchedule.Columns("N").ColumnWidth = $wsSchedule.Columns("M").ColumnWidth

# Move the selection on the Repayment schedule sheet and make it the active sheet/tab
$wsSchedule.Activate()
$wsSchedule.Range("K17").Select()

# The Transactions sheet is no longer the selected tab (its selection stays on D4)
$wsTransactions = $wb.Worksheets.Item("Transactions")
$wsTransactions.Range("D4").Select()

# Re-activate Repayment schedule so it ends up as the workbook's active tab
$wsSchedule.Activate()
